$wb = $excel.ActiveWorkbook

# ALC row 4
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 369.9375
$ws.Cells.Item(4, 9).Value = 279.85715
$ws.Cells.Item(4, 11).Value = 279.85715
$ws.Cells.Item(4, 13).Value = -165.85715

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 1269.9305
$ws.Cells.Item(137, 9).Value = 1111.7333
$ws.Cells.Item(137, 10).Value = 1533.5927
$ws.Cells.Item(137, 11).Value = 3335.199900000001
$ws.Cells.Item(137, 12).Value = 4600.7781
$ws.Cells.Item(137, 13).Value = -785.1999000000005
$ws.Cells.Item(137, 14).Value = -9700.7781

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 682.28125
$ws.Cells.Item(2, 9).Value = 531.6799999999999
$ws.Cells.Item(2, 10).Value = 1220.1428
$ws.Cells.Item(2, 11).Value = 531.6799999999999
$ws.Cells.Item(2, 12).Value = 1220.1428
$ws.Cells.Item(2, 13).Value = -418.6799999999999
$ws.Cells.Item(2, 14).Value = -1446.1428

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4290.47
$ws.Cells.Item(32, 9).Value = 3689.682
$ws.Cells.Item(32, 10).Value = 8696.25
$ws.Cells.Item(32, 11).Value = 3689.682
$ws.Cells.Item(32, 12).Value = 8696.25
$ws.Cells.Item(32, 13).Value = -3402.682
$ws.Cells.Item(32, 14).Value = -9270.25

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 83334710
$ws.Cells.Item(61, 9).Value = 111112140
$ws.Cells.Item(61, 11).Value = 111112140
$ws.Cells.Item(61, 13).Value = -111111928

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(116, 8).Value = 682.28125
$ws.Cells.Item(116, 9).Value = 531.6799999999999
$ws.Cells.Item(116, 10).Value = 1220.1428
$ws.Cells.Item(116, 11).Value = 531.6799999999999
$ws.Cells.Item(116, 12).Value = 1220.1428
$ws.Cells.Item(116, 13).Value = 1762.32
$ws.Cells.Item(116, 14).Value = -5808.1428

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 2971.7778
$ws.Cells.Item(132, 9).Value = 2559
$ws.Cells.Item(132, 11).Value = 7677
$ws.Cells.Item(132, 13).Value = -5147

# ARM row 133
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(133, 8).Value = 34120
$ws.Cells.Item(133, 10).Value = 34120
$ws.Cells.Item(133, 12).Value = 34120
$ws.Cells.Item(133, 14).Value = -39180

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 83334710
$ws.Cells.Item(136, 9).Value = 111112140
$ws.Cells.Item(136, 11).Value = 333336420
$ws.Cells.Item(136, 13).Value = -333333870

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 682.28125
$ws.Cells.Item(3, 9).Value = 531.6799999999999
$ws.Cells.Item(3, 10).Value = 1220.1428
$ws.Cells.Item(3, 11).Value = 531.6799999999999
$ws.Cells.Item(3, 12).Value = 1220.1428
$ws.Cells.Item(3, 13).Value = -417.6799999999999
$ws.Cells.Item(3, 14).Value = -1448.1428

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2406.0833
$ws.Cells.Item(20, 9).Value = 2443.9092
$ws.Cells.Item(20, 11).Value = 2443.9092
$ws.Cells.Item(20, 13).Value = -2196.9092

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 1119.3684
$ws.Cells.Item(107, 9).Value = 1069.375
$ws.Cells.Item(107, 10).Value = 1155.7273
$ws.Cells.Item(107, 11).Value = 1069.375
$ws.Cells.Item(107, 12).Value = 1155.7273
$ws.Cells.Item(107, 13).Value = 850.625
$ws.Cells.Item(107, 14).Value = -4995.7273

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1324.409
$ws.Cells.Item(31, 9).Value = 1324.409
$ws.Cells.Item(31, 10).Value = 0
$ws.Cells.Item(31, 11).Value = 1324.409
$ws.Cells.Item(31, 12).Value = 0
$ws.Cells.Item(31, 13).ClearContents()
$ws.Cells.Item(31, 14).Value = -1029.409

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 1324.409
$ws.Cells.Item(34, 9).Value = 1324.409
$ws.Cells.Item(34, 10).Value = 0
$ws.Cells.Item(34, 11).Value = 1324.409
$ws.Cells.Item(34, 12).Value = 0
$ws.Cells.Item(34, 13).ClearContents()
$ws.Cells.Item(34, 14).Value = -1122.409

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 5927.385
$ws.Cells.Item(58, 9).Value = 1136.5
$ws.Cells.Item(58, 11).Value = 1136.5
$ws.Cells.Item(58, 13).Value = -933.5

# CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 18187686
$ws.Cells.Item(62, 9).Value = 9225
$ws.Cells.Item(62, 11).Value = 9225
$ws.Cells.Item(62, 13).Value = -8601

# CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(65, 8).Value = 18187686
$ws.Cells.Item(65, 9).Value = 9225
$ws.Cells.Item(65, 11).Value = 46125
$ws.Cells.Item(65, 13).Value = -43005

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 1537.8125
$ws.Cells.Item(99, 9).Value = 1433.75
$ws.Cells.Item(99, 11).Value = 1433.75
$ws.Cells.Item(99, 13).Value = 64.25

# CRP row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(105, 8).Value = 1650
$ws.Cells.Item(105, 9).Value = 1800
$ws.Cells.Item(105, 11).Value = 1800
$ws.Cells.Item(105, 13).Value = -53

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(107, 8).Value = 970.5
$ws.Cells.Item(107, 9).Value = 502.54544
$ws.Cells.Item(107, 10).Value = 2000
$ws.Cells.Item(107, 11).Value = 502.54544
$ws.Cells.Item(107, 12).Value = 2000
$ws.Cells.Item(107, 13).Value = 1417.45456
$ws.Cells.Item(107, 14).Value = -5840

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(126, 8).Value = 1537.8125
$ws.Cells.Item(126, 9).Value = 1433.75
$ws.Cells.Item(126, 11).Value = 4301.25
$ws.Cells.Item(126, 13).Value = -1831.25

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 2383.35
$ws.Cells.Item(132, 9).Value = 1908.091
$ws.Cells.Item(132, 10).Value = 2964.2222
$ws.Cells.Item(132, 11).Value = 5724.272999999999
$ws.Cells.Item(132, 12).Value = 8892.6666
$ws.Cells.Item(132, 13).Value = -3194.272999999999
$ws.Cells.Item(132, 14).Value = -13952.6666

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 50002668
$ws.Cells.Item(134, 9).Value = 3492.75
$ws.Cells.Item(134, 10).Value = 83335450
$ws.Cells.Item(134, 11).Value = 10478.25
$ws.Cells.Item(134, 12).Value = 250006350
$ws.Cells.Item(134, 13).Value = -7943.25
$ws.Cells.Item(134, 14).Value = -250011420

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 5927.385
$ws.Cells.Item(136, 9).Value = 1136.5
$ws.Cells.Item(136, 11).Value = 3409.5
$ws.Cells.Item(136, 13).Value = -859.5

# CUL row 3
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 15053.125
$ws.Cells.Item(3, 10).Value = 19009
$ws.Cells.Item(3, 12).Value = 57027
$ws.Cells.Item(3, 14).Value = -57251

# CUL row 12
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 116.666664
$ws.Cells.Item(12, 10).Value = 72.38461
$ws.Cells.Item(12, 12).Value = 217.15383
$ws.Cells.Item(12, 14).Value = -563.15383

# CUL row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 4394.32
$ws.Cells.Item(107, 10).Value = 5646.263
$ws.Cells.Item(107, 12).Value = 16938.789
$ws.Cells.Item(107, 14).Value = -20778.789

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 45003710
$ws.Cells.Item(70, 10).Value = 40003576
$ws.Cells.Item(70, 12).Value = 40003576
$ws.Cells.Item(70, 14).Value = -40004116

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(73, 8).Value = 45003710
$ws.Cells.Item(73, 10).Value = 40003576
$ws.Cells.Item(73, 12).Value = 40003576
$ws.Cells.Item(73, 14).Value = -40005448

# GSM row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 220.76923
$ws.Cells.Item(107, 9).Value = 228.66667
$ws.Cells.Item(107, 11).Value = 228.66667
$ws.Cells.Item(107, 13).Value = 1691.33333

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 3856.5
$ws.Cells.Item(122, 9).Value = 3977.8
$ws.Cells.Item(122, 11).Value = 11933.4
$ws.Cells.Item(122, 13).Value = -9483.400000000001

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 2269.9285
$ws.Cells.Item(126, 9).Value = 1842.7142
$ws.Cells.Item(126, 11).Value = 5528.142599999999
$ws.Cells.Item(126, 13).Value = -3058.142599999999

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 3253.8147
$ws.Cells.Item(132, 9).Value = 3742.3845
$ws.Cells.Item(132, 10).Value = 2800.1428
$ws.Cells.Item(132, 11).Value = 11227.1535
$ws.Cells.Item(132, 12).Value = 8400.428400000001
$ws.Cells.Item(132, 13).Value = -8697.1535
$ws.Cells.Item(132, 14).Value = -13460.4284

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2595.6924
$ws.Cells.Item(7, 9).Value = 2529.25
$ws.Cells.Item(7, 10).Value = 2702
$ws.Cells.Item(7, 11).Value = 2529.25
$ws.Cells.Item(7, 12).Value = 2702
$ws.Cells.Item(7, 13).Value = -2417.25
$ws.Cells.Item(7, 14).Value = -2926

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 41668500
$ws.Cells.Item(122, 9).Value = 125001000
$ws.Cells.Item(122, 10).Value = 2251
$ws.Cells.Item(122, 11).Value = 375003000
$ws.Cells.Item(122, 12).Value = 6753
$ws.Cells.Item(122, 13).Value = -375000550
$ws.Cells.Item(122, 14).Value = -11653

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(126, 8).Value = 2595.6924
$ws.Cells.Item(126, 9).Value = 2529.25
$ws.Cells.Item(126, 10).Value = 2702
$ws.Cells.Item(126, 11).Value = 7587.75
$ws.Cells.Item(126, 12).Value = 8106
$ws.Cells.Item(126, 13).Value = -5117.75
$ws.Cells.Item(126, 14).Value = -13046

# LTW row 133
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(133, 8).Value = 46037
$ws.Cells.Item(133, 10).Value = 46037
$ws.Cells.Item(133, 12).Value = 46037
$ws.Cells.Item(133, 14).Value = -51097

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 2343.889
$ws.Cells.Item(136, 9).Value = 2099.1667
$ws.Cells.Item(136, 11).Value = 6297.500100000001
$ws.Cells.Item(136, 13).Value = -3747.500100000001

# WVR row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 369.3684
$ws.Cells.Item(113, 9).Value = 255.53847
$ws.Cells.Item(113, 10).Value = 616
$ws.Cells.Item(113, 11).Value = 766.61541
$ws.Cells.Item(113, 12).Value = 1848
$ws.Cells.Item(113, 13).Value = 1403.38459
$ws.Cells.Item(113, 14).Value = -6188

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 15627242
$ws.Cells.Item(122, 9).Value = 17859476
$ws.Cells.Item(122, 11).Value = 53578428
$ws.Cells.Item(122, 13).Value = -53575978

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 2332.9092
$ws.Cells.Item(132, 9).Value = 2154.7585
$ws.Cells.Item(132, 10).Value = 3624.5
$ws.Cells.Item(132, 11).Value = 6464.2755
$ws.Cells.Item(132, 12).Value = 10873.5
$ws.Cells.Item(132, 13).Value = -3934.2755
$ws.Cells.Item(132, 14).Value = -15933.5

# WVR row 133
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(133, 8).Value = 32158.75
$ws.Cells.Item(133, 10).Value = 32158.75
$ws.Cells.Item(133, 12).Value = 32158.75
$ws.Cells.Item(133, 14).Value = -42278.75

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 1707.8462
$ws.Cells.Item(136, 9).Value = 1425.875
$ws.Cells.Item(136, 10).Value = 2159
$ws.Cells.Item(136, 11).Value = 4277.625
$ws.Cells.Item(136, 12).Value = 6477
$ws.Cells.Item(136, 13).Value = -1727.625
$ws.Cells.Item(136, 14).Value = -11577
